# Append two new node rows to the schema sheet (rows 71 and 72).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 71: "a new one"
$ws.Cells.Item(71, 1).Value = 70
$ws.Cells.Item(71, 2).Value = "a new one"
$ws.Cells.Item(71, 3).Value = "a new one"
$ws.Cells.Item(71, 4).Value = 6
$ws.Cells.Item(71, 5).Value = "Custom"
$ws.Cells.Item(71, 6).Value = "a_new_one"
$ws.Cells.Item(71, 7).Value = "custom"
$ws.Cells.Item(71, 8).Value = "template"
$ws.Cells.Item(71, 9).Value = "[]"
$ws.Cells.Item(71, 10).Value = ""
$ws.Cells.Item(71, 11).Value = "['node']"
$ws.Cells.Item(71, 12).Value = "template/"

# Row 72: "a brand new node"
$ws.Cells.Item(72, 1).Value = 71
$ws.Cells.Item(72, 2).Value = "a brand new node"
$ws.Cells.Item(72, 3).Value = "a brand new node description"
$ws.Cells.Item(72, 4).Value = 6
$ws.Cells.Item(72, 5).Value = "Custom"
$ws.Cells.Item(72, 6).Value = "a_brand_new_node"
$ws.Cells.Item(72, 7).Value = "custom"
$ws.Cells.Item(72, 8).Value = "template"
$ws.Cells.Item(72, 9).Value = "[{'penalty': 'l2'}, {'C': 1.0}]"
$ws.Cells.Item(72, 10).Value = ""
$ws.Cells.Item(72, 11).Value = "['node']"
$ws.Cells.Item(72, 12).Value = "template/"
